$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row labels (old suffix -> FV2304, new suffix -> FV2310)
$old = @("Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old","Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old")
$new = @("Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new","Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new")

for ($i = 0; $i -lt $old.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $old[$i].Replace("_old", "_FV2304")
}

for ($i = 0; $i -lt $new.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $new[$i].Replace("_new", "_FV2310")
}
